# Applies the "Added timers, comments and stats" edit to Nomic.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Players" - clear out example/sample data (the "ZgH" row + values)
# ---------------------------------------------------------------------------
$players = $wb.Worksheets.Item("Players")
$players.Range("A1").Value = "Player Display Name"
$players.Range("B1").Value = ""

$players.Range("A2").Value = "Player Username"
$players.Range("B2").Value = ""

$players.Range("A3").Value = "Player ID"
$players.Range("B3").Value = ""

$players.Range("A5").Value = "Still Playing?"
$players.Range("B5").Value = ""

$players.Range("A7").Value = "Stats"

$players.Range("A8").Value = "Total Messages"
$players.Range("B8").Value = ""

$players.Range("A9").Value = "Total Days Playing"
$players.Range("B9").Value = ""

$players.Range("A10").Value = "Total Days Online"
$players.Range("B10").Value = ""

$players.Range("A11").Value = "Total Proposals"
$players.Range("B11").Value = ""

$players.Range("C20").Select()

# ---------------------------------------------------------------------------
# Sheet "Turns" - clear out example/sample data
# ---------------------------------------------------------------------------
$turns = $wb.Worksheets.Item("Turns")
$turns.Range("A1").Value = "Player Display Name"
$turns.Range("B1").Value = ""

$turns.Range("A3").Value = "Turn"
$turns.Range("C3").Value = "Proponent ID"
$turns.Range("D3").Value = "Proponent Display Name"
$turns.Range("E3").Value = "Passed?"

$turns.Range("A4").Value = ""
$turns.Range("B4").Value = ""
$turns.Range("C4").Value = ""
$turns.Range("D4").Value = ""
$turns.Range("E4").Value = ""

$turns.Range("A5").Value = ""
$turns.Range("C5").Value = ""
$turns.Range("D5").Value = ""

$turns.Range("B26").Select()

# ---------------------------------------------------------------------------
# Sheet "Misc" - add timer/stat config rows
# ---------------------------------------------------------------------------
$misc = $wb.Worksheets.Item("Misc")
$misc.Range("A1").Value = "Player Number"
$misc.Range("B1").Value = 0

$misc.Range("A3").Value = "Turn"
$misc.Range("B3").Value = 0

$misc.Range("A4").Value = "Global Turn"
$misc.Range("B4").Value = 0

$misc.Range("A5").Value = "State"
$misc.Range("B5").Value = 0

$misc.Range("A7").Value = "Current Logins"

$misc.Range("A8").Value = "Current Votes"
$misc.Range("B8").Value = 0

$misc.Range("A9").Value = "Empty?"
$misc.Range("B9").Value = 0

$misc.Range("A11").Value = "Stat Rows"
$misc.Range("B11").Value = 4

$misc.Range("A13").Value = "Proposal Time"
$misc.Range("B13").Value = 86400

$misc.Range("A14").Value = "Voting Time"
$misc.Range("B14").Value = 86400

$misc.Range("A16").Value = "Yes Proportion"
$misc.Range("B16").Value = 1

$misc.Range("B7").Select()
